# Update countries & provincias Spain
# Applies the 12-May-2020 00:35 refresh of the "Pais" COVID-19 dashboard sheet:
#   - bumps the "Datos actualizados..." timestamp from 00:05 to 00:35
#   - refreshes the case counters for several rows
#   - three countries (Colombia, Maldivas, Congo) jumped ahead of their
#     neighbours in the "Casos totales" sort, so the rows they vacated
#     shift down by one and keep their own (unchanged) figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 00:35"

# --- Helper: write a full data row (country + 7 numeric columns) ------
function Set-Row($r, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 1).Value = $country
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# --- Rows updated in place (no re-sorting needed) ----------------------
Set-Row 4  "Estados Unidos"  1381982 14344 260280 1040136 16451 779 81566
Set-Row 10 "Alemania"        172576  697   145600 19315   1576  92  7661
Set-Row 16 "Canada"          69958   1110  32934  32032   502   122 4992
Set-Row 20 "Arabia Saudita"  41014   1966  12737  28022   149   9   255

# --- Colombia overtakes Filipinas --------------------------------------
Set-Row 41 "Colombia"  11613 550 2825 8309 130 16 479
Set-Row 42 "Filipinas" 11086 292 1999 8361 31  7  726

# --- Maldivas overtakes Albania and Sri Lanka --------------------------
Set-Row 103 "Maldivas"   897 62 29  865 2 0 3
Set-Row 104 "Albania"    872 4  654 187 7 0 31
Set-Row 105 "Sri Lanka"  863 7  343 511 1 0 9

# --- Congo overtakes Mauricio, Isla de Man, Montenegro, Republica del
#     Chad, Benin, Vietnam and Ruanda -----------------------------------
Set-Row 131 "Congo"                333 59 53  269 0  1 11
Set-Row 132 "Mauricio"             332 0  322 0   0  0 10
Set-Row 133 "Isla de Man"          330 0  271 36  21 0 23
Set-Row 134 "Montenegro"           324 0  294 21  2  0 9
Set-Row 135 "Republica del Chad"   322 0  53  238 0  0 31
Set-Row 136 "Benin"                319 0  62  255 0  0 2
Set-Row 137 "Vietnam"              288 0  249 39  2  0 0
Set-Row 138 "Ruanda"               285 1  150 135 0  0 0
